$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell, resetting the cell's formatting to the
# plain/default style first (mirrors the "typed fresh data into a previously
# blank, specially-styled cell" behaviour seen in the target workbook, where the
# newly filled cells fall back to the default style instead of inheriting the
# template formatting that was sitting in the blank row).
function Set-PlainText($addr, $text) {
    $ws.Range($addr).Clear()
    $ws.Range($addr).Value = $text
}

# Helper: write a numeric-looking value as genuine TEXT (not a number), while
# still leaving the cell on the default/general style - matches phone numbers
# stored as shared-string text in the target.
function Set-PlainTextNumber($addr, $text) {
    $ws.Range($addr).Clear()
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).NumberFormat = "General"
}

# Row 5 - new entry: David / GG / PO Box 87301, Park Place, Houston, Texas
$ws.Range("B5").Value = "8/16/2023"
Set-PlainText "C5" "David"
Set-PlainText "D5" "GG"
Set-PlainText "E5" "PO Box 87301, Park Place, Houston, Texas"
Set-PlainText "F5" "PO Box 87301, Park Place, Houston, Texas"
Set-PlainTextNumber "G5" "12814104622"

# Row 6 - new entry: Henry Chelegbor / GG / 6613 Guyer Street, Philadelphia, PA, Pennsylvania
$ws.Range("B6").Value = "8/16/2023"
Set-PlainText "C6" "Henry Chelegbor"
Set-PlainText "D6" "GG"
Set-PlainText "E6" "6613 Guyer Street, Philadelphia, PA, Pennsylvania"
Set-PlainText "F6" "6613 Guyer Street, Philadelphia, PA, Pennsylvania"
Set-PlainTextNumber "G6" "13024705411"

# Row 7 - new entry: Dennis Vanmeter / GG / 1909 Harper Rd, Beckley, WV 25801
$ws.Range("B7").Value = "8/16/2023"
Set-PlainText "C7" "Dennis Vanmeter"
Set-PlainText "D7" "GG"
Set-PlainText "E7" "1909 Harper Rd, Beckley, WV 25801"
Set-PlainText "F7" "1909 Harper Rd, Beckley, WV 25801"
Set-PlainTextNumber "G7" "13049196111"

# Move the active selection to B7, matching the saved view state.
$ws.Range("B7").Select()
